# FMOD events and some Unity audio
# Created all the fmod events and started implementing some of the sound effects in Unity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (E) for rows that moved from "Not Done" to "Made"
$ws.Range("E11").Value = "Made"   # Bear Growl
$ws.Range("E12").Value = "Made"   # Bear Roar
$ws.Range("E25").Value = "Made"   # Fall No Damage

# Update Status column (E) for the row that moved from "Made" to "Done"
$ws.Range("E31").Value = "Done"   # Jumping

# Update the view: scroll so column B is the leftmost visible column,
# and select E31 (matches author's final cursor position)
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E31").Select()
